# Scheduled runner update: refresh computed market-profit values in the
# per-job Leve tracking sheets (currentAveragePrice*/LevePrice*/LeveProfit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 408814.7
$ws.Range("I28").Value = 594022.5600000001
$ws.Range("J28").Value = 15247.875
$ws.Range("K28").Value = 594022.5600000001
$ws.Range("L28").Value = 15247.875
$ws.Range("M28").Value = -593537.5600000001
$ws.Range("N28").Value = -16217.875
$ws.Range("H43").Value = 1890.8823
$ws.Range("J43").Value = 1865.8
$ws.Range("L43").Value = 1865.8
$ws.Range("N43").Value = -2003.8
$ws.Range("H70").Value = 1198
$ws.Range("I70").Value = 699
$ws.Range("J70").Value = 1697
$ws.Range("K70").Value = 2097
$ws.Range("L70").Value = 5091
$ws.Range("M70").Value = -1827
$ws.Range("N70").Value = -5631
$ws.Range("H73").Value = 1198
$ws.Range("I73").Value = 699
$ws.Range("J73").Value = 1697
$ws.Range("K73").Value = 2097
$ws.Range("L73").Value = 5091
$ws.Range("M73").Value = -1161
$ws.Range("N73").Value = -6963
$ws.Range("H80").Value = 541.5
$ws.Range("I80").Value = 604
$ws.Range("K80").Value = 1812
$ws.Range("M80").Value = -814
$ws.Range("H83").Value = 541.5
$ws.Range("I83").Value = 604
$ws.Range("K83").Value = 5436
$ws.Range("M83").Value = -444
$ws.Range("H92").Value = 1382.4
$ws.Range("I92").Value = 1228.5834
$ws.Range("J92").Value = 1997.6666
$ws.Range("K92").Value = 1228.5834
$ws.Range("L92").Value = 1997.6666
$ws.Range("M92").Value = 19.41660000000002
$ws.Range("N92").Value = -4493.6666
$ws.Range("H127").Value = 4728.2173
$ws.Range("I127").Value = 4728.2173
$ws.Range("K127").Value = 14184.6519
$ws.Range("M127").Value = -9224.651900000001
$ws.Range("H132").Value = 279885.84
$ws.Range("I132").Value = 372126.8
$ws.Range("K132").Value = 1116380.4
$ws.Range("M132").Value = -1113850.4
$ws.Range("H138").Value = 3583.4946
$ws.Range("J138").Value = 4023.9734
$ws.Range("L138").Value = 12071.9202
$ws.Range("N138").Value = -22351.9202
$ws.Range("H141").Value = 1471.4546
$ws.Range("I141").Value = 1354.1111
$ws.Range("J141").Value = 1999.5
$ws.Range("K141").Value = 4062.3333
$ws.Range("L141").Value = 5998.5
$ws.Range("M141").Value = 1117.6667
$ws.Range("N141").Value = -16358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1503020.9
$ws.Range("I32").Value = 758.0635
$ws.Range("K32").Value = 758.0635
$ws.Range("M32").Value = -471.0635
$ws.Range("H61").Value = 4088.6858
$ws.Range("I61").Value = 4088.6858
$ws.Range("K61").Value = 4088.6858
$ws.Range("M61").Value = -3876.6858
$ws.Range("H102").Value = 3477
$ws.Range("I102").Value = 3375.4666
$ws.Range("K102").Value = 3375.4666
$ws.Range("M102").Value = -1753.4666
$ws.Range("H136").Value = 4088.6858
$ws.Range("I136").Value = 4088.6858
$ws.Range("K136").Value = 12266.0574
$ws.Range("M136").Value = -9716.057400000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 526.7143
$ws.Range("I22").Value = 331.5
$ws.Range("J22").Value = 1698
$ws.Range("K22").Value = 331.5
$ws.Range("L22").Value = 1698
$ws.Range("M22").Value = -158.5
$ws.Range("N22").Value = -2044
$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -20778
$ws.Range("H99").Value = 7606.456
$ws.Range("I99").Value = 7939.881
$ws.Range("J99").Value = 6672.8667
$ws.Range("K99").Value = 7939.881
$ws.Range("L99").Value = 6672.8667
$ws.Range("M99").Value = -6441.881
$ws.Range("N99").Value = -9668.866699999999
$ws.Range("H133").Value = 91938.2
$ws.Range("J133").Value = 92672.75
$ws.Range("L133").Value = 92672.75
$ws.Range("N133").Value = -102792.75
$ws.Range("H134").Value = 921951.0600000001
$ws.Range("I134").Value = 976968.9399999999
$ws.Range("K134").Value = 2930906.82
$ws.Range("M134").Value = -2928371.82

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2994.1428
$ws.Range("J4").Value = 2994.1428
$ws.Range("L4").Value = 2994.1428
$ws.Range("N4").Value = -3218.1428
$ws.Range("H58").Value = 90922950
$ws.Range("I58").Value = 142866350
$ws.Range("J58").Value = 22020
$ws.Range("K58").Value = 142866350
$ws.Range("L58").Value = 22020
$ws.Range("M58").Value = -142866147
$ws.Range("N58").Value = -22426
$ws.Range("H62").Value = 9095.799999999999
$ws.Range("I62").Value = 4053.5
$ws.Range("K62").Value = 4053.5
$ws.Range("M62").Value = -3429.5
$ws.Range("H65").Value = 9095.799999999999
$ws.Range("I65").Value = 4053.5
$ws.Range("K65").Value = 20267.5
$ws.Range("M65").Value = -17147.5
$ws.Range("H132").Value = 7818.3887
$ws.Range("I132").Value = 8518.071
$ws.Range("J132").Value = 5369.5
$ws.Range("K132").Value = 25554.213
$ws.Range("L132").Value = 16108.5
$ws.Range("M132").Value = -23024.213
$ws.Range("N132").Value = -21168.5
$ws.Range("H135").Value = 79559.664
$ws.Range("J135").Value = 79559.664
$ws.Range("L135").Value = 79559.664
$ws.Range("N135").Value = -89699.664
$ws.Range("H136").Value = 90922950
$ws.Range("I136").Value = 142866350
$ws.Range("J136").Value = 22020
$ws.Range("K136").Value = 428599050
$ws.Range("L136").Value = 66060
$ws.Range("M136").Value = -428596500
$ws.Range("N136").Value = -71160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1134.3334
$ws.Range("I93").Value = 1097.2222
$ws.Range("J93").Value = 1178.8667
$ws.Range("K93").Value = 1097.2222
$ws.Range("L93").Value = 1227.0714
$ws.Range("M93").Value = 150.7778000000001
$ws.Range("N93").Value = -3674.8667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I2").Value = 30049.5
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 30049.5
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -29937.5
$ws.Range("N2").Value = -1224
$ws.Range("H62").Value = 18016
$ws.Range("I62").Value = 9571.6
$ws.Range("K62").Value = 9571.6
$ws.Range("M62").Value = -8947.6
$ws.Range("H65").Value = 18016
$ws.Range("I65").Value = 9571.6
$ws.Range("K65").Value = 47858
$ws.Range("M65").Value = -44738
$ws.Range("H96").Value = 1944.5
$ws.Range("I96").Value = 899
$ws.Range("K96").Value = 899
$ws.Range("M96").Value = 474
$ws.Range("H100").Value = 945.55554
$ws.Range("I100").Value = 460.8
$ws.Range("K100").Value = 921.6
$ws.Range("M100").Value = -380.6
